$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header formatting (bold font + border + centered alignment)
# from H1 onto the two new header cells so I1/J1 share the same style
# index as the rest of row 1.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats

# New header labels
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# New data values in row 2
$ws.Range("I2").Value = 8
$ws.Range("J2").Value = 8
